$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in the Price column whose new value is being written. These were
# stored as plain text in the source data (even when they look numeric, e.g.
# "211.45"), so force a text number format before writing the value to stop
# Excel from auto-converting them into floating point numbers.
$priceCells = @(
    'D2',
    'D3',
    'D5',
    'D8',
    'D9',
    'D12',
    'D13',
    'D14',
    'D15',
    'D16',
    'D18',
    'D19',
    'D20',
    'D23',
    'D24',
    'D26',
    'D27',
    'D28',
    'D31',
    'D32',
    'D34',
    'D35',
    'D38',
    'D39',
    'D43',
    'D44',
    'D45',
    'D47',
    'D48',
    'D49',
    'D50',
    'D51',
)
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '28.511.74'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '1.563.76'
$ws.Range('E3').Value = '  -1.74%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = '211.45'
$ws.Range('E5').Value = '  -1.54%  '
$ws.Range('E6').Value = '  -1.01%  '
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('D8').Value = '46.36'
$ws.Range('E8').Value = '  +5.47%  '
$ws.Range('D9').Value = '24.09'
$ws.Range('E9').Value = '  +0.27%  '
$ws.Range('E10').Value = '  -1.85%  '
$ws.Range('E11').Value = '  -1.54%  '
$ws.Range('D12').Value = '0.0887'
$ws.Range('E12').Value = '  -0.18%  '
$ws.Range('D13').Value = '1.788.52'
$ws.Range('E13').Value = '  -1.65%  '
$ws.Range('D14').Value = '1.584.02'
$ws.Range('E14').Value = '  -0.48%  '
$ws.Range('D15').Value = '0.520'
$ws.Range('E15').Value = '  -2.67%  '
$ws.Range('D16').Value = '28.524.86'
$ws.Range('E16').Value = '  +0.24%  '
$ws.Range('E17').Value = '  -3.37%  '
$ws.Range('D18').Value = '62.10'
$ws.Range('E18').Value = '  -1.73%  '
$ws.Range('D19').Value = '228.87'
$ws.Range('E19').Value = '  -1.73%  '
$ws.Range('D20').Value = '0.0₃0693'
$ws.Range('E20').Value = '  -2.42%  '
$ws.Range('E21').Value = '  -2.47%  '
$ws.Range('E22').Value = '  +0.11%  '
$ws.Range('D23').Value = '3.86'
$ws.Range('E23').Value = '  -6.37%  '
$ws.Range('D24').Value = '9.15'
$ws.Range('E24').Value = '  -3.21%  '
$ws.Range('E25').Value = '  +7.10%  '
$ws.Range('D26').Value = '150.10'
$ws.Range('E26').Value = '  -1.44%  '
$ws.Range('D27').Value = '14.96'
$ws.Range('E27').Value = '  -2.34%  '
$ws.Range('D28').Value = '6.43'
$ws.Range('E28').Value = '  -2.74%  '
$ws.Range('E29').Value = '  -4.02%  '
$ws.Range('E30').Value = '  +0.18%  '
$ws.Range('D31').Value = '0.0465'
$ws.Range('E31').Value = '  -2.04%  '
$ws.Range('D32').Value = '1.11'
$ws.Range('E32').Value = '  -3.65%  '
$ws.Range('E33').Value = '  -1.43%  '
$ws.Range('D34').Value = '3.10'
$ws.Range('E34').Value = '  -2.33%  '
$ws.Range('D35').Value = '1.395.05'
$ws.Range('E35').Value = '  -1.84%  '
$ws.Range('E36').Value = '  -1.29%  '
$ws.Range('E37').Value = '  -3.41%  '
$ws.Range('D38').Value = '2.35'
$ws.Range('E38').Value = '  +0.96%  '
$ws.Range('D39').Value = '2.59'
$ws.Range('E39').Value = '  +2.36%  '
$ws.Range('E40').Value = '  -1.20%  '
$ws.Range('E41').Value = '  -1.60%  '
$ws.Range('E42').Value = '  +0.11%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = '1.88'
$ws.Range('E43').Value = '  +2.43%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').Value = '0.788'
$ws.Range('E44').Value = '  -4.11%  '
$ws.Range('D45').Value = '5.51'
$ws.Range('E45').Value = '  -4.56%  '
$ws.Range('E46').Value = '  -0.02%  '
$ws.Range('D47').Value = '62.69'
$ws.Range('E47').Value = '  -3.21%  '
$ws.Range('D48').Value = '1.700.23'
$ws.Range('E48').Value = '  -1.78%  '
$ws.Range('D49').Value = '86.13'
$ws.Range('E49').Value = '  -1.85%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.0522'
$ws.Range('E50').Value = '  -0.25%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0102'
$ws.Range('E51').Value = '  -5.58%  '

# Restore the default (unstyled) appearance so only the cell content changes,
# matching the original workbook where these cells had no explicit style.
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).Style = "Normal"
}
